# Auto-generated Excel COM-interop script to apply the Leve profit data refresh
# (chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 165.83333
$ws.Range("I4").Value = 165.83333
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 165.83333
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -51.83332999999999
$ws.Range("N4").Value = ""
# Row 41
$ws.Range("H41").Value = 339.5
$ws.Range("J41").Value = 373.35715
$ws.Range("L41").Value = 373.35715
$ws.Range("N41").Value = -1253.35715
# Row 57
$ws.Range("H57").Value = 22156
$ws.Range("J57").Value = 22156
$ws.Range("L57").Value = 66468
$ws.Range("N57").Value = -67466
# Row 62
$ws.Range("H62").Value = 2273.3684
$ws.Range("I62").Value = 1799.6154
$ws.Range("K62").Value = 1799.6154
$ws.Range("M62").Value = -1175.6154
# Row 65
$ws.Range("H65").Value = 2273.3684
$ws.Range("I65").Value = 1799.6154
$ws.Range("K65").Value = 8998.076999999999
$ws.Range("M65").Value = -5878.076999999999
# Row 86
$ws.Range("H86").Value = 6559.7896
$ws.Range("I86").Value = 1189.2858
$ws.Range("K86").Value = 1189.2858
$ws.Range("M86").Value = -66.28580000000011
# Row 89
$ws.Range("H89").Value = 6559.7896
$ws.Range("I89").Value = 1189.2858
$ws.Range("K89").Value = 5946.429
$ws.Range("M89").Value = -330.4290000000001
# Row 100
$ws.Range("H100").Value = 2275.5557
$ws.Range("I100").Value = 1640
$ws.Range("J100").Value = 4500
$ws.Range("K100").Value = 1640
$ws.Range("L100").Value = 4500
$ws.Range("M100").Value = -1099
$ws.Range("N100").Value = -5582
# Row 129
$ws.Range("H129").Value = 847.5571
$ws.Range("J129").Value = 847.4853000000001
$ws.Range("L129").Value = 2542.4559
$ws.Range("N129").Value = -12542.4559

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3533.0967
$ws.Range("I45").Value = 4031.2
$ws.Range("J45").Value = 3295.9048
$ws.Range("K45").Value = 4031.2
$ws.Range("L45").Value = 3295.9048
$ws.Range("M45").Value = -3654.2
$ws.Range("N45").Value = -4049.9048
# Row 61
$ws.Range("H61").Value = 859964.5600000001
$ws.Range("I61").Value = 1287588.8
$ws.Range("J61").Value = 4716.2856
$ws.Range("K61").Value = 1287588.8
$ws.Range("L61").Value = 4716.2856
$ws.Range("M61").Value = -1287376.8
$ws.Range("N61").Value = -5140.2856
# Row 74
$ws.Range("H74").Value = 1756.9
$ws.Range("I74").Value = 1818.7778
$ws.Range("K74").Value = 1818.7778
$ws.Range("M74").Value = -944.7778000000001
# Row 77
$ws.Range("H77").Value = 1756.9
$ws.Range("I77").Value = 1818.7778
$ws.Range("K77").Value = 9093.889000000001
$ws.Range("M77").Value = -4725.889000000001
# Row 102
$ws.Range("H102").Value = 3184.5386
$ws.Range("I102").Value = 1403.3334
$ws.Range("J102").Value = 4711.2856
$ws.Range("K102").Value = 1403.3334
$ws.Range("L102").Value = 4711.2856
$ws.Range("M102").Value = 218.6666
$ws.Range("N102").Value = -7955.2856
# Row 108
$ws.Range("H108").Value = 32000
$ws.Range("J108").Value = 32000
$ws.Range("L108").Value = 32000
$ws.Range("N108").Value = -39680
# Row 122
$ws.Range("H122").Value = 2001.5
$ws.Range("I122").Value = 2073.1428
$ws.Range("K122").Value = 6219.428400000001
$ws.Range("M122").Value = -3769.428400000001
# Row 136
$ws.Range("H136").Value = 859964.5600000001
$ws.Range("I136").Value = 1287588.8
$ws.Range("J136").Value = 4716.2856
$ws.Range("K136").Value = 3862766.4
$ws.Range("L136").Value = 14148.8568
$ws.Range("M136").Value = -3860216.4
$ws.Range("N136").Value = -19248.8568

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 10
$ws.Range("H10").Value = 100
$ws.Range("I10").Value = 100
$ws.Range("K10").Value = 100
$ws.Range("M10").Value = 40
# Row 105
$ws.Range("H105").Value = 2098.4
$ws.Range("I105").Value = 2098.5
$ws.Range("J105").Value = 2098.2856
$ws.Range("K105").Value = 2098.5
$ws.Range("L105").Value = 2098.2856
$ws.Range("M105").Value = -351.5
$ws.Range("N105").Value = -5592.2856

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 12
$ws.Range("H12").Value = 377.5
$ws.Range("I12").Value = 255
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 255
$ws.Range("L12").Value = 500
$ws.Range("M12").Value = -85
$ws.Range("N12").Value = -840
# Row 31
$ws.Range("H31").Value = 9798.195
$ws.Range("I31").Value = 11799.257
$ws.Range("K31").Value = 11799.257
$ws.Range("M31").Value = -11504.257
# Row 34
$ws.Range("H34").Value = 9798.195
$ws.Range("I34").Value = 11799.257
$ws.Range("K34").Value = 11799.257
$ws.Range("M34").Value = -11597.257
# Row 47
$ws.Range("H47").Value = 10933.333
$ws.Range("I47").Value = 5000
$ws.Range("J47").Value = 13900
$ws.Range("K47").Value = 5000
$ws.Range("L47").Value = 13900
$ws.Range("M47").Value = -4434
$ws.Range("N47").Value = -15032
# Row 93
$ws.Range("H93").Value = 7576.75
$ws.Range("J93").Value = 12000
$ws.Range("L93").Value = 12000
$ws.Range("N93").Value = -15744
# Row 99
$ws.Range("H99").Value = 4522.8696
$ws.Range("I99").Value = 3324.2354
$ws.Range("K99").Value = 3324.2354
$ws.Range("M99").Value = -1826.2354
# Row 103
$ws.Range("H103").Value = 18353.4
$ws.Range("I103").Value = 7643
$ws.Range("J103").Value = 25493.666
$ws.Range("K103").Value = 7643
$ws.Range("L103").Value = 25493.666
$ws.Range("M103").Value = -6471
$ws.Range("N103").Value = -27837.666
# Row 122
$ws.Range("H122").Value = 1533.2222
$ws.Range("J122").Value = 1525
$ws.Range("L122").Value = 4575
$ws.Range("N122").Value = -9475
# Row 126
$ws.Range("H126").Value = 4522.8696
$ws.Range("I126").Value = 3324.2354
$ws.Range("K126").Value = 9972.706200000001
$ws.Range("M126").Value = -7502.706200000001
# Row 134
$ws.Range("H134").Value = 1292.2727
$ws.Range("I134").Value = 1038.6154
$ws.Range("J134").Value = 1658.6666
$ws.Range("K134").Value = 3115.8462
$ws.Range("L134").Value = 4975.9998
$ws.Range("M134").Value = -580.8462
$ws.Range("N134").Value = -10045.9998
# Row 135
$ws.Range("H135").Value = 50516.668
$ws.Range("J135").Value = 50516.668
$ws.Range("L135").Value = 50516.668
$ws.Range("N135").Value = -60656.668

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 109
$ws.Range("H109").Value = 2369.7144
$ws.Range("I109").Value = 1369
$ws.Range("K109").Value = 4107
$ws.Range("M109").Value = -3067
# Row 113
$ws.Range("H113").Value = 502.1111
$ws.Range("I113").Value = 519
$ws.Range("J113").Value = 493.66666
$ws.Range("K113").Value = 1557
$ws.Range("L113").Value = 1480.99998
$ws.Range("M113").Value = 613
$ws.Range("N113").Value = -5820.999980000001
# Row 131
$ws.Range("H131").Value = 727.12
$ws.Range("J131").Value = 749.49475
$ws.Range("L131").Value = 2248.48425
$ws.Range("N131").Value = -12328.48425

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 71.333336
$ws.Range("I2").Value = 45.833332
$ws.Range("J2").Value = 173.33333
$ws.Range("K2").Value = 45.833332
$ws.Range("L2").Value = 173.33333
$ws.Range("M2").Value = 67.166668
$ws.Range("N2").Value = -399.33333
# Row 17
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").Value = ""
# Row 96
$ws.Range("H96").Value = 19630.5
$ws.Range("J96").Value = 19630.5
$ws.Range("L96").Value = 19630.5
$ws.Range("N96").Value = -25122.5
# Row 122
$ws.Range("H122").Value = 1727.8889
$ws.Range("I122").Value = 1413.75
$ws.Range("K122").Value = 4241.25
$ws.Range("M122").Value = -1791.25
# Row 132
$ws.Range("H132").Value = 35779.375
$ws.Range("I132").Value = 5838.857
$ws.Range("J132").Value = 59066.445
$ws.Range("K132").Value = 17516.571
$ws.Range("L132").Value = 177199.335
$ws.Range("M132").Value = -14986.571
$ws.Range("N132").Value = -182259.335

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3739.524
$ws.Range("I40").Value = 3208.125
$ws.Range("J40").Value = 5440
$ws.Range("K40").Value = 3208.125
$ws.Range("L40").Value = 5440
$ws.Range("M40").Value = -3072.125
$ws.Range("N40").Value = -5712
# Row 51
$ws.Range("H51").Value = 19000
$ws.Range("J51").Value = 19000
$ws.Range("L51").Value = 19000
$ws.Range("N51").Value = -19956
# Row 100
$ws.Range("H100").Value = 2667.0908
$ws.Range("I100").Value = 1790.4615
$ws.Range("J100").Value = 3933.3333
$ws.Range("K100").Value = 1790.4615
$ws.Range("L100").Value = 3933.3333
$ws.Range("M100").Value = -1249.4615
$ws.Range("N100").Value = -5015.3333
# Row 109
$ws.Range("H109").Value = 31661.666
$ws.Range("J109").Value = 31661.666
$ws.Range("L109").Value = 31661.666
$ws.Range("N109").Value = -34435.666
# Row 122
$ws.Range("H122").Value = 983952.4
$ws.Range("J122").Value = 3830.7693
$ws.Range("L122").Value = 11492.3079
$ws.Range("N122").Value = -16392.3079
# Row 136
$ws.Range("H136").Value = 118538.38
$ws.Range("I136").Value = 168555.44
$ws.Range("K136").Value = 505666.32
$ws.Range("M136").Value = -503116.32

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 61
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 19000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 19000
$ws.Range("M61").Value = ""
$ws.Range("N61").Value = -19584
# Row 92
$ws.Range("H92").Value = 17207.143
$ws.Range("J92").Value = 17207.143
$ws.Range("L92").Value = 17207.143
$ws.Range("N92").Value = -22199.143
# Row 94
$ws.Range("H94").Value = 31990
$ws.Range("J94").Value = 31990
$ws.Range("L94").Value = 31990
$ws.Range("N94").Value = -33792
# Row 97
$ws.Range("H97").Value = 43000
$ws.Range("J97").Value = 43000
$ws.Range("L97").Value = 43000
$ws.Range("N97").Value = -44982
# Row 126
$ws.Range("H126").Value = 1385.1333
$ws.Range("I126").Value = 848.375
$ws.Range("J126").Value = 1998.5714
$ws.Range("K126").Value = 2545.125
$ws.Range("L126").Value = 5995.7142
$ws.Range("M126").Value = -75.125
$ws.Range("N126").Value = -10935.7142
# Row 136
$ws.Range("H136").Value = 1204.3948
$ws.Range("I136").Value = 911.0417
$ws.Range("J136").Value = 1707.2858
$ws.Range("K136").Value = 2733.1251
$ws.Range("L136").Value = 5121.857400000001
$ws.Range("M136").Value = -183.1251000000002
$ws.Range("N136").Value = -10221.8574

Write-Output "Applied Leve profit updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets"